$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension -> measure renames (and two cells whose value changes to a
# brand-new "residencia-*-nombre" measure string)
$ws.Range("D2").Value = "iaest-measure:residencia-comarca-nombre"
$ws.Range("E2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:residencia-provincia-nombre"
$ws.Range("H2").Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Range("I2").Value = "iaest-measure:edad-grandes-grupos"
$ws.Range("J2").Value = "iaest-measure:relacion-lugar-de-residencia-y-nacimiento"

# Row 3: these columns switch from "dim" to "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: these columns switch from "skos:Concept" / "URI-comarca" / "URI-Provincia" to "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: the per-dimension mapping file references for the now-curated
# dimensions/measures are no longer needed, so clear those cells entirely
# (Clear, not ClearContents, so the cell element itself is dropped).
$ws.Range("E5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
$ws.Range("J5").Clear()
